$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Настройки и Условия")

# Rewrite the "preparation conditions" block (rows 33-38) with the new,
# shorter wording. Row 39 ("Товар со скидкой") is removed/merged in, so
# that text now lives on row 37 and the old row 39 cell becomes empty.
$ws.Range("A33").Value = "Товар без наличия"
$ws.Range("A34").Value = "Товар без наличия и с действием “Предзаказ”"
$ws.Range("A35").Value = "Товар без цены"
$ws.Range("A36").Value = "Товар без цены и с действием “Попросить покупателя ввести цену”"
$ws.Range("A37").Value = "Товар со скидкой"
$ws.Range("A38").Value = "Подготовить бренд — вывод логотипа на стр категории"

# Row 39 no longer holds content - clear it back to a plain, unformatted
# empty cell (also resets its style to the default).
$ws.Range("A39").ClearContents()
$ws.Range("A39").NumberFormat = "General"

# The active selection on this sheet moved to A39.
$ws.Range("A39").Select()
